$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the promo text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.0 = 7251.51 pesos`n✅ 7251.51 pesos = 1.98 = 822.69 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 500.999
$ws2.Range("O10").Value = 3633
$ws2.Range("N12").Value = 3658
$ws2.Range("O12").Value = 415.001
